# Auto-generated edit script applying the Tonberry_Profits.xlsx diff
# Updates computed market-board profit columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 879688.8
$ws.Range("I92").Value = 1026214.06
$ws.Range("J92").Value = 537
$ws.Range("K92").Value = 1026214.06
$ws.Range("L92").Value = 537
$ws.Range("M92").Value = -1024966.06
$ws.Range("N92").Value = -3033

$ws.Range("H100").Value = 1845.9
$ws.Range("I100").Value = 1717.6666
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1717.6666
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1176.6666
$ws.Range("N100").Value = -4082


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 382.85715
$ws.Range("I5").Value = 222
$ws.Range("K5").Value = 222
$ws.Range("M5").Value = -110

$ws.Range("H32").Value = 3706.1904
$ws.Range("I32").Value = 3041.0205
$ws.Range("K32").Value = 3041.0205
$ws.Range("M32").Value = -2754.0205

$ws.Range("H45").Value = 1580.3572
$ws.Range("I45").Value = 1050.4
$ws.Range("J45").Value = 1874.7778
$ws.Range("K45").Value = 1050.4
$ws.Range("L45").Value = 1874.7778
$ws.Range("M45").Value = -673.4000000000001
$ws.Range("N45").Value = -2628.7778

$ws.Range("H61").Value = 4359.5557
$ws.Range("I61").Value = 3192.3333
$ws.Range("J61").Value = 6694
$ws.Range("K61").Value = 3192.3333
$ws.Range("L61").Value = 6694
$ws.Range("M61").Value = -2980.3333
$ws.Range("N61").Value = -7118

$ws.Range("H74").Value = 1010.44684
$ws.Range("I74").Value = 850.7317
$ws.Range("K74").Value = 850.7317
$ws.Range("M74").Value = 23.26829999999995

$ws.Range("H77").Value = 1010.44684
$ws.Range("I77").Value = 850.7317
$ws.Range("K77").Value = 4253.6585
$ws.Range("M77").Value = 114.3414999999995

$ws.Range("H97").Value = 760.8461
$ws.Range("I97").Value = 599
$ws.Range("K97").Value = 599
$ws.Range("M97").Value = -103

$ws.Range("H122").Value = 1590.5834
$ws.Range("I122").Value = 1553.3636
$ws.Range("K122").Value = 4660.0908
$ws.Range("M122").Value = -2210.0908

$ws.Range("H132").Value = 1988.8928
$ws.Range("I132").Value = 1592.7333
$ws.Range("K132").Value = 4778.199900000001
$ws.Range("M132").Value = -2248.199900000001

$ws.Range("H136").Value = 4359.5557
$ws.Range("I136").Value = 3192.3333
$ws.Range("J136").Value = 6694
$ws.Range("K136").Value = 9576.999899999999
$ws.Range("L136").Value = 20082
$ws.Range("M136").Value = -7026.999899999999
$ws.Range("N136").Value = -25182


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 382.85715
$ws.Range("I4").Value = 222
$ws.Range("K4").Value = 222
$ws.Range("M4").Value = -107

$ws.Range("H20").Value = 2587.3333
$ws.Range("I20").Value = 2335.875
$ws.Range("J20").Value = 4599
$ws.Range("K20").Value = 2335.875
$ws.Range("L20").Value = 4599
$ws.Range("M20").Value = -2088.875
$ws.Range("N20").Value = -5093

$ws.Range("H86").Value = 97686.52
$ws.Range("I86").Value = 2247.8572
$ws.Range("J86").Value = 288563.84
$ws.Range("K86").Value = 2247.8572
$ws.Range("L86").Value = 288563.84
$ws.Range("M86").Value = -1124.8572
$ws.Range("N86").Value = -290809.84

$ws.Range("H89").Value = 97686.52
$ws.Range("I89").Value = 2247.8572
$ws.Range("J89").Value = 288563.84
$ws.Range("K89").Value = 11239.286
$ws.Range("L89").Value = 1442819.2
$ws.Range("M89").Value = -5623.286
$ws.Range("N89").Value = -1454051.2

$ws.Range("H134").Value = 6448.25
$ws.Range("I134").Value = 7189.625
$ws.Range("K134").Value = 21568.875
$ws.Range("M134").Value = -19033.875


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1892577.5
$ws.Range("I58").Value = 2900210.5
$ws.Range("J58").Value = 3265.375
$ws.Range("K58").Value = 2900210.5
$ws.Range("L58").Value = 3265.375
$ws.Range("M58").Value = -2900007.5
$ws.Range("N58").Value = -3671.375

$ws.Range("H132").Value = 2468.182
$ws.Range("I132").Value = 1090.75
$ws.Range("J132").Value = 4121.1
$ws.Range("K132").Value = 3272.25
$ws.Range("L132").Value = 12363.3
$ws.Range("M132").Value = -742.25
$ws.Range("N132").Value = -17423.3

$ws.Range("H134").Value = 1442.1578
$ws.Range("I134").Value = 1455.6111
$ws.Range("K134").Value = 4366.8333
$ws.Range("M134").Value = -1831.8333

$ws.Range("H136").Value = 1892577.5
$ws.Range("I136").Value = 2900210.5
$ws.Range("J136").Value = 3265.375
$ws.Range("K136").Value = 8700631.5
$ws.Range("L136").Value = 9796.125
$ws.Range("M136").Value = -8698081.5
$ws.Range("N136").Value = -14896.125


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 159.83333
$ws.Range("I6").Value = 71.8
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 215.4
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -102.4
$ws.Range("N6").Value = -2026

$ws.Range("H11").Value = 1229.3334
$ws.Range("I11").Value = 1200
$ws.Range("J11").Value = 1244
$ws.Range("K11").Value = 3600
$ws.Range("L11").Value = 3732
$ws.Range("M11").Value = -3460
$ws.Range("N11").Value = -4012

$ws.Range("H21").Value = 1166.6666
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1166.6666
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3499.9998
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -3845.9998

$ws.Range("H56").Value = 10734
$ws.Range("I56").Value = 10734
$ws.Range("K56").Value = 10734
$ws.Range("M56").Value = -10204

$ws.Range("H129").Value = 26714.678
$ws.Range("I129").Value = 420.55554
$ws.Range("J129").Value = 39169.79
$ws.Range("K129").Value = 1261.66662
$ws.Range("L129").Value = 117509.37
$ws.Range("M129").Value = 3738.33338
$ws.Range("N129").Value = -127509.37

$ws.Range("H131").Value = 782.97
$ws.Range("J131").Value = 793.71875
$ws.Range("L131").Value = 2381.15625
$ws.Range("N131").Value = -12461.15625

$ws.Range("H132").Value = 1712.25
$ws.Range("I132").Value = 699.5
$ws.Range("K132").Value = 6295.5
$ws.Range("M132").Value = -3765.5

$ws.Range("H140").Value = 1493.3793
$ws.Range("I140").Value = 865.8421
$ws.Range("J140").Value = 2685.7
$ws.Range("K140").Value = 2597.5263
$ws.Range("L140").Value = 8057.099999999999
$ws.Range("M140").Value = 2582.4737
$ws.Range("N140").Value = -18417.1


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 117.29412
$ws.Range("I2").Value = 152
$ws.Range("J2").Value = 86.44444
$ws.Range("K2").Value = 152
$ws.Range("L2").Value = 86.44444
$ws.Range("M2").Value = -39
$ws.Range("N2").Value = -312.44444

$ws.Range("H80").Value = 3347.8
$ws.Range("I80").Value = 2916.6667
$ws.Range("K80").Value = 2916.6667
$ws.Range("M80").Value = -1918.6667

$ws.Range("H83").Value = 3347.8
$ws.Range("I83").Value = 2916.6667
$ws.Range("K83").Value = 14583.3335
$ws.Range("M83").Value = -9591.333500000001

$ws.Range("H132").Value = 3499174.2
$ws.Range("I132").Value = 6412123.5
$ws.Range("J132").Value = 3635.4
$ws.Range("K132").Value = 19236370.5
$ws.Range("L132").Value = 10906.2
$ws.Range("M132").Value = -19233840.5
$ws.Range("N132").Value = -15966.2

$ws.Range("H138").Value = 42332.25
$ws.Range("I138").Value = 25150
$ws.Range("K138").Value = 25150
$ws.Range("M138").Value = -20010


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2202.1667
$ws.Range("I22").Value = 2566
$ws.Range("J22").Value = 1942.2858
$ws.Range("K22").Value = 2566
$ws.Range("L22").Value = 1942.2858
$ws.Range("M22").Value = -2271
$ws.Range("N22").Value = -2532.2858

$ws.Range("H27").Value = 2202.1667
$ws.Range("I27").Value = 2566
$ws.Range("J27").Value = 1942.2858
$ws.Range("K27").Value = 2566
$ws.Range("L27").Value = 1942.2858
$ws.Range("M27").Value = -2459
$ws.Range("N27").Value = -2156.2858

$ws.Range("H43").Value = 12008.75
$ws.Range("J43").Value = 12008.75
$ws.Range("L43").Value = 12008.75
$ws.Range("N43").Value = -12394.75

$ws.Range("H58").Value = 10697.667
$ws.Range("I58").Value = 3093
$ws.Range("J58").Value = 14500
$ws.Range("K58").Value = 3093
$ws.Range("L58").Value = 14500
$ws.Range("M58").Value = -2833
$ws.Range("N58").Value = -15020

$ws.Range("H68").Value = 2774
$ws.Range("J68").Value = 3233
$ws.Range("L68").Value = 3233
$ws.Range("N68").Value = -4731

$ws.Range("H71").Value = 2774
$ws.Range("J71").Value = 3233
$ws.Range("L71").Value = 16165
$ws.Range("N71").Value = -23653

$ws.Range("H82").Value = 1883.8334
$ws.Range("I82").Value = 1303.5
$ws.Range("K82").Value = 1303.5
$ws.Range("M82").Value = -942.5

$ws.Range("H85").Value = 1883.8334
$ws.Range("I85").Value = 1303.5
$ws.Range("K85").Value = 1303.5
$ws.Range("M85").Value = -55.5

$ws.Range("H132").Value = 1646.025
$ws.Range("I132").Value = 1321
$ws.Range("J132").Value = 1911.9546
$ws.Range("K132").Value = 3963
$ws.Range("L132").Value = 5735.8638
$ws.Range("M132").Value = -1433
$ws.Range("N132").Value = -10795.8638


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 61499.383
$ws.Range("I122").Value = 112309.71
$ws.Range("J122").Value = 2220.6667
$ws.Range("K122").Value = 336929.13
$ws.Range("L122").Value = 6662.000100000001
$ws.Range("M122").Value = -334479.13
$ws.Range("N122").Value = -11562.0001

$ws.Range("H130").Value = 34580.832
$ws.Range("J130").Value = 34580.832
$ws.Range("L130").Value = 34580.832
$ws.Range("N130").Value = -44620.832

$ws.Range("H132").Value = 1497.5264
$ws.Range("I132").Value = 1173.7646
$ws.Range("K132").Value = 3521.2938
$ws.Range("M132").Value = -991.2937999999999

$ws.Range("H136").Value = 2585.96
$ws.Range("I136").Value = 2796.5
$ws.Range("J136").Value = 2318
$ws.Range("K136").Value = 8389.5
$ws.Range("L136").Value = 6954
$ws.Range("M136").Value = -5839.5
$ws.Range("N136").Value = -12054

